# Append a new data row (row 16) to the IBBRandom sheet, matching the
# existing table layout: Date | totalScore | ... | Method

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16

# Column A: date/time value (stored as Excel serial date, same style as A2:A15)
$ws.Cells.Item($row, 1).Value = 42625.886689814812

# Column B: totalScore
$ws.Cells.Item($row, 2).Value = 23

# Columns C..M: the various percentage/count metrics, all zero for this row
for ($col = 3; $col -le 13; $col++) {
    $ws.Cells.Item($row, $col).Value = 0
}

# Column N: Method label
$ws.Cells.Item($row, 14).Value = "Random"
